# Generate Report for Handoff
#
# A new localization handoff run completed, producing a fresh content GUID
# and target-file hash. Refresh the three report sheets (Overview, zh-cn,
# de-de) so the "latest handoff" file names / timestamps reflect this run.

$wb = $excel.ActiveWorkbook

$newGuid = "13155f2b-d060-4c87-b808-502c3c655232"
$newHash = "97e03b65603dc46efb9910ff43750fe743d81ff1"

# --- Overview sheet: source file name + latest handoff date ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("D2").Value = "2016-50-12 20:50:13"

# --- zh-cn sheet: source file name, handoff target file, handoff datetime ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("D2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-12 20:50:09"

# --- de-de sheet: source file name, handoff target file, handoff datetime ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("D2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-12 20:50:13"
